# Apply BOM V2.1 updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: capacitor designator list gains C75 and C67
$ws.Range("B2").Value = "C1,C36,C53,C4,C18,C6,C17,C75,C43,C62,C24,C28,C38,C61,C52,C50,C47,C67,C3,C8,C27,C34,C58,C29"

# Row 8: designator list C68,C67 -> C68,C74
$ws.Range("B8").Value = "C68,C74"

# Row 10: capacitor designator list gains C73
$ws.Range("B10").Value = "C70,C41,C55,C45,C66,C46,C49,C31,C2,C22,C54,C59,C73,C56,C65,C64,C9,C71,C33,C57,C16,C72,C63,C23"

# Row 15: designator list J5,J4,J3 -> J4,J3 (J5 removed)
$ws.Range("B15").Value = "J4,J3"

# Row 16: connector footprint/part changed, LCSC part number added
$ws.Range("C16").Value = "TF3822S05SV830"
$ws.Range("D16").Value = "C6552695"

# Row 29: new component added (IMU)
$ws.Range("A29").Value = "ICM-42688-P"
$ws.Range("B29").Value = "U8"
$ws.Range("C29").Value = "PQFN50P300X250X97-14N"
$ws.Range("D29").Value = "C1850418"

$wb.Save()
